# ---------------------------------------------------------------------------
# Add 2022-Q3 data:
#  1. Insert a new row into the "总计" (summary) sheet for the 2022-Q3 figures
#     and renumber the existing rows' index column.
#  2. Insert a brand-new worksheet named "2022-Q3" right after "总计",
#     populated with the Q3 fund-holdings table (copying the layout/format
#     from the neighbouring "2022-Q2" sheet).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make room for the new 2022-Q3 row right under the header row.
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

# Re-apply the same formatting used by the other index cells in column A.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 2.47

# Renumber the index column for the rows that shifted down by one.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $summary)
$newSheet.Name = "2022-Q3"

# Reuse the layout/formatting of the "2022-Q2" sheet as a starting point.
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Range("A1:H10").Copy($newSheet.Range("A1"))
$newSheet.Range("A8:H10").Clear()
$newSheet.Range("A1").Clear()

# Row 2
$newSheet.Range("B2").Value = "'501208"
$newSheet.Range("C2").Value = "中欧创新未来混合（LOF）"
$newSheet.Range("D2").Value = "'55.17"
$newSheet.Range("E2").Value = "'85.30"
$newSheet.Range("F2").Value = "'3.22"
$newSheet.Range("G2").Value = "'1.7765"
$newSheet.Range("H2").Value = 8

# Row 3
$newSheet.Range("B3").Value = "'501081"
$newSheet.Range("C3").Value = "中欧科创主题混合（LOF）"
$newSheet.Range("D3").Value = "'7.39"
$newSheet.Range("E3").Value = "'86.44"
$newSheet.Range("F3").Value = "'4.77"
$newSheet.Range("G3").Value = "'0.3525"
$newSheet.Range("H3").Value = 6

# Row 4
$newSheet.Range("B4").Value = "'610001"
$newSheet.Range("C4").Value = "信澳领先增长混合A"
$newSheet.Range("D4").Value = "'7.32"
$newSheet.Range("E4").Value = "'92.87"
$newSheet.Range("F4").Value = "'3.74"
$newSheet.Range("G4").Value = "'0.2738"
$newSheet.Range("H4").Value = 7

# Row 5
$newSheet.Range("B5").Value = "'015143"
$newSheet.Range("C5").Value = "中欧智能制造混合A"
$newSheet.Range("D5").Value = "'0.96"
$newSheet.Range("E5").Value = "'84.58"
$newSheet.Range("F5").Value = "'3.46"
$newSheet.Range("G5").Value = "'0.0332"
$newSheet.Range("H5").Value = 6

# Row 6
$newSheet.Range("B6").Value = "'015456"
$newSheet.Range("C6").Value = "信澳领先增长混合C"
$newSheet.Range("D6").Value = "'0.55"
$newSheet.Range("E6").Value = "'92.87"
$newSheet.Range("F6").Value = "'3.74"
$newSheet.Range("G6").Value = "'0.0206"
$newSheet.Range("H6").Value = 7

# Row 7
$newSheet.Range("B7").Value = "'015144"
$newSheet.Range("C7").Value = "中欧智能制造混合C"
$newSheet.Range("D7").Value = "'0.51"
$newSheet.Range("E7").Value = "'84.58"
$newSheet.Range("F7").Value = "'3.46"
$newSheet.Range("G7").Value = "'0.0176"
$newSheet.Range("H7").Value = 6

# The apostrophe-prefixed assignments above forced plain numeric-looking
# strings to be stored as text, but that also stamps a "quote prefix" style
# onto the cell. These data cells carry no explicit styling in the source
# table, so strip that incidental formatting back off (this does not affect
# the stored text values).
$newSheet.Range("B2:B7").ClearFormats()
$newSheet.Range("D2:G7").ClearFormats()
